$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "Research Questions" slide at position 2 (right after
#    the title slide). Layout 16 == "Title and Content" (matches the layout
#    used by every other content slide in this deck).
# ---------------------------------------------------------------------------
$rq = $p.Slides.Add(2, 16)

$rqTitle = $rq.Shapes.Item(1).TextFrame.TextRange
$rqTitle.Text = "Research Questions"
$rqTitle.LanguageID = "en-AU"

$rqBody = $rq.Shapes.Item(2).TextFrame.TextRange
$rqBody.Text = "MAIN Question: how did Australia break the world record for economic growth`rSupplementary Question: what are the main factors that contribute to economic growth in Australia`rSupplementary Question:"
$rqBody.LanguageID = "en-AU"

# Paragraph 1 - bold the "MAIN Question:" lead-in only.
$rqP1 = $rqBody.Paragraphs(1, 1)
$rqP1Lead = $rqP1.Characters(1, "MAIN Question:".Length)
$rqP1Lead.Font.Bold = $true

# Paragraph 2 - demote to level 2, bold+italic the "Supplementary Question: " lead-in.
$rqP2 = $rqBody.Paragraphs(2, 1)
$rqP2.IndentLevel = 2
$rqP2Lead = $rqP2.Characters(1, "Supplementary Question: ".Length)
$rqP2Lead.Font.Bold = $true
$rqP2Lead.Font.Italic = $true

# Paragraph 3 - demote to level 2, entire line bold+italic.
$rqP3 = $rqBody.Paragraphs(3, 1)
$rqP3.IndentLevel = 2
$rqP3.Font.Bold = $true
$rqP3.Font.Italic = $true

# ---------------------------------------------------------------------------
# 2. The old slide 2 ("How do you measure the economy") is now slide 3 because
#    of the insertion above. Append three new bullet paragraphs to its body
#    content placeholder, right before the existing trailing blank paragraph.
# ---------------------------------------------------------------------------
$gdp = $p.Slides.Item(3)
$gdpBody = $gdp.Shapes.Item(2).TextFrame.TextRange

$lastBullet = $gdpBody.Paragraphs($gdpBody.Paragraphs().Count, 1)
$lastBullet.InsertAfter("`rGDP is measured on a quarterly basis (or every 3 months)`rEconomic growth occurs when GDP increases between quarters (e.g. Q2 > Q1)`rAustralia has avoided a recession (defined as two consecutive quarters of negative GDP growth) for 28 years to hold the world record for continued economic growth")
